$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The MARCA (E) and MODELO (F) columns held brand names that are no longer
# supplied by the model; VALOR UNITARIO (G) was a placeholder price. The
# new algorithm that locates the respective MARCA/MODELO/VALOR columns
# found nothing to fill them with, so clear these three columns for every
# data row (rows 2 through 110).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 7).Value = ""
}
